$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (leading space is intentional, matches original naming)
$ws.Name = " FT_L7901"

# Move the active selection to A5
$ws.Range("A5").Select()
